$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-02-16 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-02-17 Saturday", 2)

$d.Content.Find.Execute("33×53=1749", $true, $false, $false, $false, $false, $true, 1, $false, "31×55=1705", 2)
$d.Content.Find.Execute("85×88=7480", $true, $false, $false, $false, $false, $true, 1, $false, "48×61=2928", 2)
$d.Content.Find.Execute("79×87=6873", $true, $false, $false, $false, $false, $true, 1, $false, "84×65=5460", 2)
$d.Content.Find.Execute("73×18=1314", $true, $false, $false, $false, $false, $true, 1, $false, "27×49=1323", 2)
$d.Content.Find.Execute("21×65=1365", $true, $false, $false, $false, $false, $true, 1, $false, "30×32=960", 2)

$d.Content.Find.Execute("60×86=5160", $true, $false, $false, $false, $false, $true, 1, $false, "75×19=1425", 2)
$d.Content.Find.Execute("35×89=3115", $true, $false, $false, $false, $false, $true, 1, $false, "41×70=2870", 2)
$d.Content.Find.Execute("68×47=3196", $true, $false, $false, $false, $false, $true, 1, $false, "71×73=5183", 2)
$d.Content.Find.Execute("17×81=1377", $true, $false, $false, $false, $false, $true, 1, $false, "71×55=3905", 2)
$d.Content.Find.Execute("34×67=2278", $true, $false, $false, $false, $false, $true, 1, $false, "27×99=2673", 2)

$d.Content.Find.Execute("12×65=780", $true, $false, $false, $false, $false, $true, 1, $false, "16×82=1312", 2)
$d.Content.Find.Execute("75×43=3225", $true, $false, $false, $false, $false, $true, 1, $false, "74×45=3330", 2)
$d.Content.Find.Execute("14×57=798", $true, $false, $false, $false, $false, $true, 1, $false, "11×13=143", 2)
$d.Content.Find.Execute("52×71=3692", $true, $false, $false, $false, $false, $true, 1, $false, "54×60=3240", 2)
$d.Content.Find.Execute("72×13=936", $true, $false, $false, $false, $false, $true, 1, $false, "50×35=1750", 2)

$d.Content.Find.Execute("97×84=8148", $true, $false, $false, $false, $false, $true, 1, $false, "77×31=2387", 2)
$d.Content.Find.Execute("44×66=2904", $true, $false, $false, $false, $false, $true, 1, $false, "54×31=1674", 2)
$d.Content.Find.Execute("63×85=5355", $true, $false, $false, $false, $false, $true, 1, $false, "56×43=2408", 2)
$d.Content.Find.Execute("82×42=3444", $true, $false, $false, $false, $false, $true, 1, $false, "41×86=3526", 2)
$d.Content.Find.Execute("39×59=2301", $true, $false, $false, $false, $false, $true, 1, $false, "60×40=2400", 2)

$d.Content.Find.Execute("49×44=2156", $true, $false, $false, $false, $false, $true, 1, $false, "68×69=4692", 2)
$d.Content.Find.Execute("54×79=4266", $true, $false, $false, $false, $false, $true, 1, $false, "50×23=1150", 2)
$d.Content.Find.Execute("15×23=345", $true, $false, $false, $false, $false, $true, 1, $false, "14×94=1316", 2)
$d.Content.Find.Execute("64×91=5824", $true, $false, $false, $false, $false, $true, 1, $false, "59×35=2065", 2)
$d.Content.Find.Execute("28×89=2492", $true, $false, $false, $false, $false, $true, 1, $false, "67×80=5360", 2)
